$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.194.16'
$ws.Range('E2').Value = '  -1.18%  '

$ws.Range('D3').Value = '1.796.19'
$ws.Range('E3').Value = '  -1.55%  '

$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.39'
$ws.Range('E5').Value = '  -0.42%  '

$ws.Range('E6').Value = '  -0.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5198'
$ws.Range('E7').Value = '  +1.84%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3817'
$ws.Range('E8').Value = '  -3.32%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07934'
$ws.Range('E9').Value = '  -3.24%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.47'
$ws.Range('E10').Value = '  -0.70%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.099'
$ws.Range('E11').Value = '  -1.24%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.282'
$ws.Range('E12').Value = '  -1.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.001'
$ws.Range('E13').Value = '  -0.14%  '

$ws.Range('E14').Value = '  -2.69%  '

$ws.Range('D15').Value = '1.796.52'
$ws.Range('E15').Value = '  -1.18%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.255'
$ws.Range('E16').Value = '  -4.07%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.19'
$ws.Range('E17').Value = '  +0.27%  '

$ws.Range('E18').Value = '  -3.02%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06555'
$ws.Range('E19').Value = '  -1.53%  '

$ws.Range('E20').Value = '  -0.05%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.31'
$ws.Range('E21').Value = '  -2.82%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.950'
$ws.Range('E22').Value = '  -2.48%  '

$ws.Range('D23').Value = '28.232.76'
$ws.Range('E23').Value = '  -1.20%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.14'
$ws.Range('E24').Value = '  -2.42%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.267'
$ws.Range('E25').Value = '  -0.10%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.97'
$ws.Range('E26').Value = '  +2.46%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.45'
$ws.Range('E27').Value = '  -4.12%  '

$ws.Range('D28').Value = '2.001.15'
$ws.Range('E28').Value = '  -1.39%  '

$ws.Range('E29').Value = '  -2.96%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.14'
$ws.Range('E30').Value = '  -2.69%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1066'
$ws.Range('E31').Value = '  -1.85%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.053'
$ws.Range('E32').Value = '  -5.61%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.670'
$ws.Range('E33').Value = '  +0.21%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.565'
$ws.Range('E34').Value = '  -3.68%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07315'
$ws.Range('E35').Value = '  +3.43%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.28'
$ws.Range('E36').Value = '  +8.85%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02327'
$ws.Range('E37').Value = '  -1.15%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2141'
$ws.Range('E38').Value = '  -4.04%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.068'
$ws.Range('E39').Value = '  -3.41%  '

$ws.Range('E40').Value = '  -1.76%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6166'
$ws.Range('E41').Value = '  -2.90%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.162'
$ws.Range('E42').Value = '  -1.63%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.368'
$ws.Range('E43').Value = '  -2.32%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.27'
$ws.Range('E44').Value = '  -1.98%  '

$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6003'
$ws.Range('E45').Value = '  +0.52%  '

$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.780'
$ws.Range('E46').Value = '  +1.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.47'
$ws.Range('E47').Value = '  +1.70%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.230'
$ws.Range('E48').Value = '  +2.96%  '

$ws.Range('E49').Value = '  -3.63%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06775'
$ws.Range('E50').Value = '  -2.43%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.13'
$ws.Range('E51').Value = '  -1.59%  '
